$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price (D) and volume-change (E) updates
$ws.Range("D2").Value = "30.050.78"
$ws.Range("E2").Value = "  -1.71%  "
$ws.Range("D3").Value = "2.102.54"
$ws.Range("E3").Value = "  -0.60%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.65%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "347.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.59%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5168"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.46%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4439"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.65%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "52.30"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.49%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08953"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.71%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.171"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.18%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "25.57"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.83%  "
$ws.Range("D13").Value = "2.111.28"
$ws.Range("E13").Value = "  -0.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.238"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.725"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "99.41"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.41%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001150"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.04%  "
$ws.Range("E18").Value = "  -0.58%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "20.78"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.94%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06683"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.10%  "
$ws.Range("E21").Value = "  -0.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.245"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.04%  "
$ws.Range("D23").Value = "30.154.30"
$ws.Range("E23").Value = "  -1.59%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.73"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.345"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.60%  "
$ws.Range("D26").Value = "2.356.40"
$ws.Range("E26").Value = "  -0.38%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.81%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.530"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.75%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "162.20"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.54"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.47%  "
$ws.Range("E31").Value = "  -3.09%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1065"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.67%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.633"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.45%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.233"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.956"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.929"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.32%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.22"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.39%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02574"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.87%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06805"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.43%  "
$ws.Range("E40").Value = "  -1.51%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.58"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6807"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.288"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.33%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.30"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.40%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6371"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.290"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.23%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000366"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.83%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.642"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.26%  "
$ws.Range("E49").Value = "  -2.78%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "82.28"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.51%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07228"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.36%  "
